$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column B ("Unnamed: 0") in front of the old "#Ronde"
#     column, duplicating the row-index values already present in column A.
#     Everything that used to be in columns B..H ("#Ronde".."EO1") shifts
#     one column to the right (C..I), and a brand new "EO1" header shows
#     up at the new last column I1.
$ws.Columns("B").Insert()

# Give the new B1 header the same look (bold / boxed) as the other header
# cells by copying the format from a neighbouring header cell, then fill
# in its text.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "Unnamed: 0"

# The new last header cell (I1) needs the same header formatting too.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "EO1"
$excel.CutCopyMode = $false

# Column B (new) mirrors column A's row index (0..19); column C keeps the
# old "#Ronde" numbers (1..20) that used to live in column B, but loses
# the old 1-off numeric formatting it had inherited from the insert/shift.
$ws.Range("B2:C21").Clear()
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 2).Value = $r - 2
    $ws.Cells.Item($r, 3).Value = $r - 1
}

$ws.Range("C27").Select()
